$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.451.43"
$ws.Range("E2").Value = "  -0.44%  "
$ws.Range("D3").Value = "1.802.25"
$ws.Range("E3").Value = "  +0.28%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "228.47"
$ws.Range("E5").Value = "  +0.66%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.580"
$ws.Range("E6").Value = "  +3.59%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "34.59"
$ws.Range("E8").Value = "  +4.81%  "
$ws.Range("E9").Value = "  +1.33%  "
$ws.Range("E10").Value = "  -0.09%  "
$ws.Range("E11").Value = "  +0.17%  "
$ws.Range("E12").Value = "  +0.29%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.21"
$ws.Range("E13").Value = "  +1.34%  "
$ws.Range("D14").Value = "1.803.38"
$ws.Range("E14").Value = "  +0.51%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.642"
$ws.Range("E15").Value = "  +0.78%  "
$ws.Range("D16").Value = "34.428.25"
$ws.Range("E16").Value = "  -0.42%  "
$ws.Range("E17").Value = "  +1.44%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "69.02"
$ws.Range("E18").Value = "  +0.10%  "
$ws.Range("D19").Value = "0.0₃0798"
$ws.Range("E19").Value = "  -0.39%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "245.20"
$ws.Range("E20").Value = "  -0.97%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.51"
$ws.Range("E21").Value = "  +1.20%  "
$ws.Range("E22").Value = "  +0.20%  "
$ws.Range("E23").Value = "  -0.19%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "174.03"
$ws.Range("E24").Value = "  +3.97%  "
$ws.Range("E25").Value = "  +1.87%  "
$ws.Range("E26").Value = "  +6.39%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.80"
$ws.Range("E27").Value = "  +1.15%  "
$ws.Range("E28").Value = "  +2.11%  "
$ws.Range("E29").Value = "  -1.02%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.01"
$ws.Range("E30").Value = "  -1.95%  "
$ws.Range("E31").Value = "  +1.00%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.24"
$ws.Range("E32").Value = "  +0.28%  "
$ws.Range("E33").Value = "  +0.49%  "
$ws.Range("E34").Value = "  -0.49%  "
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.683"
$ws.Range("E35").Value = "  +1.75%  "
$ws.Range("B36").Value = "Maker"
$ws.Range("C36").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D36").Value = "1.395.07"
$ws.Range("E36").Value = "  -2.37%  "
$ws.Range("E37").Value = "  -3.79%  "
$ws.Range("E38").Value = "  -0.38%  "
$ws.Range("E39").Value = "  -1.53%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "83.46"
$ws.Range("E41").Value = "  +2.52%  "
$ws.Range("E42").Value = "  +1.32%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.41"
$ws.Range("E43").Value = "  -0.07%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.61"
$ws.Range("E44").Value = "  -0.73%  "
$ws.Range("E45").Value = "  +3.46%  "
$ws.Range("E46").Value = "  -3.34%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "5.97"
$ws.Range("E47").Value = "  -2.39%  "
$ws.Range("D48").Value = "1.962.34"
$ws.Range("E48").Value = "  +0.33%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "104.96"
$ws.Range("E49").Value = "  -1.07%  "
$ws.Range("E50").Value = "  +0.09%  "
$ws.Range("E51").Value = "  +1.32%  "
